$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.924.22"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "2.918.23"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "597.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.550"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("E9").Value = "  +3.32%  "
$ws.Range("D10").Value = "2.917.81"
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.432"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +16.27%  "
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.87"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "3.454.79"
$ws.Range("E14").Value = "  +3.04%  "
$ws.Range("D15").Value = "75.829.69"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.36"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").Value = "2.914.14"
$ws.Range("E18").Value = "  +2.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.96"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.78"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "378.04"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.28"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "3.069.46"
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.20"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.68"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  +5.99%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "505.93"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.71"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.03"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.11"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.68"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("B39").Value = "Cronos"
$ws.Range("C39").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.106"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +23.20%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.113"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.30%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "180.42"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.343"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.98"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.66"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.21"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.573"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("E50").Value = "  +7.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.71"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.94%  "
